# Microsite Education Script completed
# Appends new interview-history rows to AMSIN, BETA and AMS sheets, and
# fixes up the timestamp/style of the existing AMS row 31 ("164_live").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write one data row (Run Date, Run Time, Sprint Name, Total,
# Pass, Fail, Time Taken) onto a worksheet at a given row number.
# Column A holds a date-looking label that must stay literal TEXT (not
# get auto-converted to a real date serial), so we force the cell to
# Text before assigning the value and then drop it back to General
# once the literal string is safely stored. Column B is the numeric
# run-time serial, formatted with the sheet's custom date-time format.
# Every cell's NumberFormat is (re)applied explicitly - including the
# plain numeric columns - so the formatting sticks even when a cell's
# value happens to be unchanged from what was already on the sheet.
# ---------------------------------------------------------------------
function Write-HistoryRow {
    param(
        $ws,
        [int]$row,
        [string]$runDate,
        [double]$runTime,
        [string]$sprintName,
        [double]$totalCases,
        [double]$passCases,
        [double]$failCases,
        [double]$timeTaken
    )

    $aCell = $ws.Cells.Item($row, 1)
    $aCell.NumberFormat = "@"
    $aCell.Value = $runDate
    $aCell.NumberFormat = "General"

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $runTime
    $bCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $sprintName
    $cCell.NumberFormat = "General"

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.Value = $totalCases
    $dCell.NumberFormat = "General"

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $passCases
    $eCell.NumberFormat = "General"

    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $failCases
    $fCell.NumberFormat = "General"

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.Value = $timeTaken
    $gCell.NumberFormat = "General"
}

# ---------------------------------------------------------------------
# AMSIN: add rows 45-49 (165 / 166 first-cycle/second-cycle/final-run
# sprints), dimension grows from A1:G44 to A1:G49.
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")
Write-HistoryRow $wsAmsin 45 "2022-08-02" 44775.64808427083 "165_fstcycle"   155 145 10 6.7
Write-HistoryRow $wsAmsin 46 "2022-08-03" 44776.66374898148 "165_scndcycle" 155 154 1  3.67
Write-HistoryRow $wsAmsin 47 "2022-08-04" 44777.38665354167 "165_finalrun"  155 152 3  3.45
Write-HistoryRow $wsAmsin 48 "2022-08-22" 44795.66772075232 "166fstcycle"   155 151 4  4.72
Write-HistoryRow $wsAmsin 49 "2022-08-23" 44796.90077827546 "166cyclescnd"  155 155 0  3.38

# ---------------------------------------------------------------------
# BETA: add rows 25-26 (165 / 166 beta runs), dimension grows from
# A1:G24 to A1:G26.
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")
Write-HistoryRow $wsBeta 25 "2022-08-04" 44777.55462648148 "165beta"  155 154 1 3.05
Write-HistoryRow $wsBeta 26 "2022-08-24" 44797.52063487269 "166_beta" 155 149 6 4.04

# ---------------------------------------------------------------------
# AMS: row 31 ("164_live") gets its run-time timestamp corrected and
# picks up explicit styling, then rows 32-33 (165 / 166 live runs) are
# appended, dimension grows from A1:G31 to A1:G33.
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")
Write-HistoryRow $wsAms 31 "2022-07-14" 44756.81479289352 "164_live" 155 153 2 3.19
Write-HistoryRow $wsAms 32 "2022-08-04" 44777.80886119213 "165_live" 155 153 2 2.91
Write-HistoryRow $wsAms 33 "2022-08-24" 44797.91623644932 "166_live" 155 149 6 4.4
